$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.923.28'
$ws.Cells.Item(2, 5).Value = '  -0.47%  '
$ws.Cells.Item(3, 4).Value = '1.879.42'
$ws.Cells.Item(3, 5).Value = '  -0.97%  '
$ws.Cells.Item(5, 4).Value = '324.83'
$ws.Cells.Item(5, 5).Value = '  -0.74%  '
$ws.Cells.Item(6, 4).Value = '1.003'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '
$ws.Cells.Item(7, 4).Value = '0.4607'
$ws.Cells.Item(7, 5).Value = '  -0.89%  '
$ws.Cells.Item(8, 4).Value = '0.3869'
$ws.Cells.Item(8, 5).Value = '  -1.26%  '
$ws.Cells.Item(9, 4).Value = '0.07853'
$ws.Cells.Item(9, 5).Value = '  -1.25%  '
$ws.Cells.Item(10, 4).Value = '0.9854'
$ws.Cells.Item(10, 5).Value = '  -2.67%  '
$ws.Cells.Item(11, 4).Value = '21.80'
$ws.Cells.Item(11, 5).Value = '  -1.44%  '
$ws.Cells.Item(12, 4).Value = '1.877.20'
$ws.Cells.Item(12, 5).Value = '  -1.65%  '
$ws.Cells.Item(13, 4).Value = '6.992'
$ws.Cells.Item(13, 5).Value = '  -1.96%  '
$ws.Cells.Item(14, 5).Value = '  -2.30%  '
$ws.Cells.Item(15, 4).Value = '0.06974'
$ws.Cells.Item(15, 5).Value = '  -0.03%  '
$ws.Cells.Item(16, 4).Value = '88.02'
$ws.Cells.Item(16, 5).Value = '  -1.35%  '
$ws.Cells.Item(17, 5).Value = '  +0.03%  '
$ws.Cells.Item(18, 4).Value = '0.000009957'
$ws.Cells.Item(18, 5).Value = '  -1.78%  '
$ws.Cells.Item(19, 4).Value = '16.93'
$ws.Cells.Item(20, 5).Value = '  -0.20%  '
$ws.Cells.Item(21, 4).Value = '28.917.28'
$ws.Cells.Item(21, 5).Value = '  -0.52%  '
$ws.Cells.Item(22, 4).Value = '5.251'
$ws.Cells.Item(22, 5).Value = '  -1.99%  '
$ws.Cells.Item(23, 4).Value = '10.96'
$ws.Cells.Item(23, 5).Value = '  -1.46%  '
$ws.Cells.Item(24, 4).Value = '2.104'
$ws.Cells.Item(24, 5).Value = '  +2.20%  '
$ws.Cells.Item(25, 4).Value = '156.34'
$ws.Cells.Item(25, 5).Value = '  +0.69%  '
$ws.Cells.Item(26, 4).Value = '19.34'
$ws.Cells.Item(26, 5).Value = '  -2.65%  '
$ws.Cells.Item(27, 4).Value = '5.985'
$ws.Cells.Item(27, 5).Value = '  +1.79%  '
$ws.Cells.Item(28, 4).Value = '117.55'
$ws.Cells.Item(28, 5).Value = '  -1.94%  '
$ws.Cells.Item(29, 4).Value = '1.913'
$ws.Cells.Item(29, 5).Value = '  -3.93%  '
$ws.Cells.Item(30, 4).Value = '0.09352'
$ws.Cells.Item(30, 5).Value = '  -0.42%  '
$ws.Cells.Item(31, 4).Value = '0.9010'
$ws.Cells.Item(31, 5).Value = '  -4.09%  '
$ws.Cells.Item(32, 4).Value = '5.265'
$ws.Cells.Item(32, 5).Value = '  -1.91%  '
$ws.Cells.Item(33, 4).Value = '1.320'
$ws.Cells.Item(33, 5).Value = '  -2.38%  '
$ws.Cells.Item(34, 5).Value = '  -0.16%  '
$ws.Cells.Item(35, 5).Value = '  +0.38%  '
$ws.Cells.Item(36, 4).Value = '0.05745'
$ws.Cells.Item(36, 5).Value = '  -1.58%  '
$ws.Cells.Item(37, 4).Value = '0.02074'
$ws.Cells.Item(37, 5).Value = '  -1.03%  '
$ws.Cells.Item(38, 5).Value = '  -0.10%  '
$ws.Cells.Item(39, 4).Value = '7.625'
$ws.Cells.Item(39, 5).Value = '  -5.92%  '
$ws.Cells.Item(40, 4).Value = '0.5645'
$ws.Cells.Item(40, 5).Value = '  -3.31%  '
$ws.Cells.Item(41, 4).Value = '0.1765'
$ws.Cells.Item(41, 5).Value = '  -2.83%  '
$ws.Cells.Item(42, 4).Value = '9.716'
$ws.Cells.Item(42, 5).Value = '  -2.84%  '
$ws.Cells.Item(43, 4).Value = '2.255'
$ws.Cells.Item(43, 5).Value = '  -0.39%  '
$ws.Cells.Item(44, 4).Value = '11.89'
$ws.Cells.Item(44, 5).Value = '  -0.74%  '
$ws.Cells.Item(45, 4).Value = '0.5340'
$ws.Cells.Item(45, 5).Value = '  -2.30%  '
$ws.Cells.Item(46, 4).Value = '0.07048'
$ws.Cells.Item(46, 5).Value = '  -2.16%  '
$ws.Cells.Item(47, 4).Value = '1.840'
$ws.Cells.Item(47, 5).Value = '  -1.60%  '
$ws.Cells.Item(48, 4).Value = '2.558'
$ws.Cells.Item(48, 5).Value = '  +2.47%  '
$ws.Cells.Item(49, 4).Value = '112.65'
$ws.Cells.Item(49, 5).Value = '  -0.70%  '
$ws.Cells.Item(50, 4).Value = '1.061'
$ws.Cells.Item(50, 5).Value = '  -5.69%  '
$ws.Cells.Item(51, 4).Value = '70.73'
$ws.Cells.Item(51, 5).Value = '  -0.78%  '
